$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.701.97"
$ws.Range("E2").Value = "  -2.70%  "
$ws.Range("D3").Value = "1.558.97"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "205.64"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "21.92"
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("D9").Value = "0.247"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").Value = "1.780.99"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "1.560.52"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").Value = "61.47"
$ws.Range("E16").Value = "  -2.96%  "
$ws.Range("D17").Value = "26.737.73"
$ws.Range("E17").Value = "  -2.46%  "
$ws.Range("D18").Value = "213.97"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").Value = "7.32"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("D24").Value = "2.00"
$ws.Range("E24").Value = "  -0.78%  "
$ws.Range("D25").Value = "152.96"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").Value = "14.83"
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("E30").Value = "  -3.85%  "
$ws.Range("E31").Value = "  -1.58%  "
$ws.Range("E32").Value = "  -1.32%  "
$ws.Range("D33").Value = "1.386.51"
$ws.Range("E33").Value = "  +1.85%  "
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("E35").Value = "  +1.55%  "
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").Value = "0.930"
$ws.Range("E37").Value = "  -4.43%  "
$ws.Range("E38").Value = "  -2.63%  "
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("D40").Value = "0.808"
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").Value = "0.991"
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("D46").Value = "63.16"
$ws.Range("E46").Value = "  -1.54%  "
$ws.Range("D47").Value = "1.694.44"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").Value = "85.56"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").Value = "0.0₇0985"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("E51").Value = "  -0.98%  "
